$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price cells in column D become plain decimal numbers after this
# update (e.g. "214.85"). Excel auto-converts numeric-looking text assigned to a
# cell, so those specific cells are pre-formatted as Text to keep them stored as
# literal strings (matching the rest of the Price column, which stays textual
# because most prices contain multiple "." separators Excel never reads as a number).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"

$ws.Range("D2").Value = '27.051.65'
$ws.Range("E2").Value = '  -0.44%  '
$ws.Range("D3").Value = '1.623.07'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '214.85'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '0.0634'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("D10").Value = '20.10'
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.850.75'
$ws.Range("D13").Value = '1.630.80'
$ws.Range("E13").Value = '  -0.49%  '
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '64.80'
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("D17").Value = '27.026.82'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '213.78'
$ws.Range("E19").Value = '  -1.62%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("E23").Value = '  -6.23%  '
$ws.Range("E24").Value = '  -0.61%  '
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -2.76%  '
$ws.Range("D29").Value = '15.60'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.762'
$ws.Range("E32").Value = '  +38.77%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '3.35'
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").Value = '1.359.94'
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("D39").Value = '0.848'
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").Value = '0.803'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.23'
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '65.37'
$ws.Range("E43").Value = '  +5.79%  '
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("D45").Value = '1.762.40'
$ws.Range("D46").Value = '0.885'
$ws.Range("E46").Value = '  +32.19%  '
$ws.Range("D47").Value = '90.15'
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  +2.86%  '
$ws.Range("D49").Value = '0.0₆0107'
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("E50").Value = '  +6.80%  '
$ws.Range("E51").Value = '  +0.48%  '
